# Insert a new weekly data record at row 58 ("Apio" / Macroferia Regional de
# Talca sheet). All existing records from row 58 down to row 131 shift down
# by one row (to rows 59-132), which is exactly what a native row insert
# does, and the sheet's used-range grows from A1:R131 to A1:R132.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 58:131 down to 59:132, leaving row 58 blank for the new record.
$ws.Rows("58:58").Insert()

# Populate the newly inserted row 58 with the new record's data.
$ws.Range("A58").Value = 5
$ws.Range("B58").Value = "Macroferia Regional de Talca"
$ws.Range("C58").Value = "Maule"
$ws.Range("D58").Value = 44494
$ws.Range("E58").Value = 7
$ws.Range("F58").Value = 100112017
$ws.Range("G58").Value = "Apio"
$ws.Range("H58").Value = "Americana (o)"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 500
$ws.Range("K58").Value = 7000
$ws.Range("L58").Value = 7000
$ws.Range("M58").Value = 7000
$ws.Range("N58").Value = "`$/docena de matas"
$ws.Range("O58").Value = "Provincia del Elquí"
$ws.Range("P58").Value = 1167
$ws.Range("Q58").Value = 6
$ws.Range("R58").Value = "Hortaliza"
